# Actualización desde MV -datos-
# Append 4 new daily rows (07-09-2021 .. 10-09-2021) to the EMBI spreads sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Serie = "07-09-2021"; Values = @(307.3, 374, 194, 300, 85.8, 114.7, 15.8, 151.7, 457, 1476.9, 282.7, 269, 133, 341, 171) },
    @{ Serie = "08-09-2021"; Values = @(310.4, 378, 195, 303, 85.7, 113.4, 16.2, 154.4, 461.6, 1504, 290, 272, 136, 342, 174) },
    @{ Serie = "09-09-2021"; Values = @(313.6, 382, 197, 306, 85.2, 114.1, 14.2, 156.7, 465.6, 1528.4, 294.1, 274, 140, 347, 175) },
    @{ Serie = "10-09-2021"; Values = @(310.4, 380, 194, 301, 84.8, 113.2, 13.8, 151.8, 460.2, 1538.2, 287.5, 273, 137, 345, 172) }
)

$startRow = 172

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    # Force the date-like label to be stored as literal text (matching the
    # source data's "dd-mm-yyyy" string series) instead of letting Excel's
    # automatic date recognition convert it to a serial date number.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $data.Serie
    $ws.Cells.Item($row, 1).Style = "Normal"

    for ($c = 0; $c -lt $data.Values.Count; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $data.Values[$c]
    }
}
